# Update odds values on existing rows 2, 3, 4, 8 and 11, and append a new
# match (row 15) at the bottom of the sheet, per the "Jogos da Semana"
# FlashScore odds refresh for 2024-11-01.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,7).Value = 2.15
$ws.Cells.Item(2,9).Value = 4
$ws.Cells.Item(2,10).Value = 3.1
$ws.Cells.Item(2,24).Value = 8
$ws.Cells.Item(2,33).Value = 7
$ws.Cells.Item(2,35).Value = 17
$ws.Cells.Item(2,36).Value = 51
$ws.Cells.Item(2,42).Value = 34
$ws.Cells.Item(2,52).Value = 126

# Row 3
$ws.Cells.Item(3,7).Value = 2.15
$ws.Cells.Item(3,17).Value = 3.1
$ws.Cells.Item(3,18).Value = 1.36
$ws.Cells.Item(3,24).Value = 8.5
$ws.Cells.Item(3,42).Value = 34
$ws.Cells.Item(3,43).Value = 51

# Row 4
$ws.Cells.Item(4,7).Value = 1.73
$ws.Cells.Item(4,21).Value = 2.5
$ws.Cells.Item(4,22).Value = 1.5
$ws.Cells.Item(4,23).Value = 4.75
$ws.Cells.Item(4,24).Value = 6.5
$ws.Cells.Item(4,35).Value = 21
$ws.Cells.Item(4,40).Value = 3.4
$ws.Cells.Item(4,43).Value = 34
$ws.Cells.Item(4,49).Value = 7
$ws.Cells.Item(4,52).Value = 151

# Row 8
$ws.Cells.Item(8,7).Value = 2.63
$ws.Cells.Item(8,9).Value = 2.8
$ws.Cells.Item(8,10).Value = 3.6
$ws.Cells.Item(8,13).Value = 1.13
$ws.Cells.Item(8,14).Value = 6
$ws.Cells.Item(8,23).Value = 6
$ws.Cells.Item(8,27).Value = 29
$ws.Cells.Item(8,45).Value = 351

# Row 11
$ws.Cells.Item(11,7).Value = 2.18
$ws.Cells.Item(11,8).Value = 3.2
$ws.Cells.Item(11,9).Value = 3.15
$ws.Cells.Item(11,10).Value = 2.8
$ws.Cells.Item(11,11).Value = 2.05
$ws.Cells.Item(11,12).Value = 3.7
$ws.Cells.Item(11,13).Value = 8.3
$ws.Cells.Item(11,14).Value = 1.05
$ws.Cells.Item(11,19).Value = 1.4
$ws.Cells.Item(11,20).Value = 2.5
$ws.Cells.Item(11,21).Value = 1.7
$ws.Cells.Item(11,22).Value = 1.93
$ws.Cells.Item(11,23).Value = 7.7
$ws.Cells.Item(11,24).Value = 10.75
$ws.Cells.Item(11,25).Value = 8.75
$ws.Cells.Item(11,27).Value = 17.5
$ws.Cells.Item(11,28).Value = 27
$ws.Cells.Item(11,31).Value = 13.5
$ws.Cells.Item(11,34).Value = 17
$ws.Cells.Item(11,35).Value = 11
$ws.Cells.Item(11,36).Value = 40
$ws.Cells.Item(11,37).Value = 27
$ws.Cells.Item(11,40).Value = 4.05
$ws.Cells.Item(11,41).Value = 11.5
$ws.Cells.Item(11,42).Value = 19.5
$ws.Cells.Item(11,44).Value = 80
$ws.Cells.Item(11,46).Value = 2.47
$ws.Cells.Item(11,49).Value = 5
$ws.Cells.Item(11,50).Value = 17.5
$ws.Cells.Item(11,51).Value = 24
$ws.Cells.Item(11,52).Value = 90
$ws.Cells.Item(11,53).Value = 120
$ws.Cells.Item(11,54).Value = 300

# Row 15 (new row)
$ws.Cells.Item(15,1).Value = "Kz2IMp1S"
$ws.Cells.Item(15,3).Value = "22:30"
$ws.Cells.Item(15,4).Value = "USA - MLS"
$ws.Cells.Item(15,5).Value = "Colorado Rapids"
$ws.Cells.Item(15,6).Value = "Los Angeles Galaxy"
$ws.Range("B15").NumberFormat = "@"
$ws.Cells.Item(15,2).Value = "01/11/2024"
$ws.Cells.Item(15,7).Value = 2.8
$ws.Cells.Item(15,8).Value = 3.9
$ws.Cells.Item(15,9).Value = 2.3
$ws.Cells.Item(15,10).Value = 3.1
$ws.Cells.Item(15,11).Value = 2.6
$ws.Cells.Item(15,12).Value = 2.75
$ws.Cells.Item(15,13).Value = 1.01
$ws.Cells.Item(15,14).Value = 23
$ws.Cells.Item(15,15).Value = 1.11
$ws.Cells.Item(15,16).Value = 6.5
$ws.Cells.Item(15,17).Value = 1.4
$ws.Cells.Item(15,18).Value = 2.88
$ws.Cells.Item(15,19).Value = 1.22
$ws.Cells.Item(15,20).Value = 4
$ws.Cells.Item(15,21).Value = 1.36
$ws.Cells.Item(15,22).Value = 3
$ws.Cells.Item(15,23).Value = 17
$ws.Cells.Item(15,24).Value = 19
$ws.Cells.Item(15,25).Value = 11
$ws.Cells.Item(15,26).Value = 29
$ws.Cells.Item(15,27).Value = 19
$ws.Cells.Item(15,28).Value = 21
$ws.Cells.Item(15,29).Value = 23
$ws.Cells.Item(15,30).Value = 8.5
$ws.Cells.Item(15,31).Value = 11
$ws.Cells.Item(15,32).Value = 26
$ws.Cells.Item(15,33).Value = 15
$ws.Cells.Item(15,34).Value = 17
$ws.Cells.Item(15,35).Value = 10
$ws.Cells.Item(15,36).Value = 23
$ws.Cells.Item(15,37).Value = 15
$ws.Cells.Item(15,38).Value = 17
$ws.Cells.Item(15,39).Value = 67
$ws.Cells.Item(15,40).Value = 5.5
$ws.Cells.Item(15,41).Value = 13
$ws.Cells.Item(15,42).Value = 17
$ws.Cells.Item(15,43).Value = 41
$ws.Cells.Item(15,44).Value = 41
$ws.Cells.Item(15,45).Value = 81
$ws.Cells.Item(15,46).Value = 4
$ws.Cells.Item(15,47).Value = 6.5
$ws.Cells.Item(15,48).Value = 34
$ws.Cells.Item(15,49).Value = 5
$ws.Cells.Item(15,50).Value = 12
$ws.Cells.Item(15,51).Value = 15
$ws.Cells.Item(15,52).Value = 34
$ws.Cells.Item(15,53).Value = 41
$ws.Cells.Item(15,54).Value = 81
$ws.Cells.Item(15,55).Value = 201
$ws.Cells.Item(15,56).Value = 176
